# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (B15:J23) gets re-sorted by column E (Periodo
# Mora) in ascending order. Columns B, C, D, G, H, I, J are identical for
# every data row, so the only cells whose *content* actually changes are
# column E (the period label) and column F (the "Valor Mora" date/serial
# value that travels together with its row).
#
# Before (rows 16-23, col E): 2107, 2106, 2105, 2104, 2103, 2102, 2101, 2012
# After  (rows 16-23, col E): 2012, 2101, 2102, 2103, 2104, 2105, 2106, 2107
#
# F16 / F23 are the only F-column values that differ between two rows in
# the original data (29260 vs 35112), so after sorting, F16 becomes 35112
# and F23 becomes 29260 while F17:F22 stay 35112.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("2012", "2101", "2102", "2103", "2104", "2105", "2106", "2107")
$valores = @(35112, 35112, 35112, 35112, 35112, 35112, 35112, 29260)

for ($i = 0; $i -lt 8; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}

# The autofit ("best fit") columns B:J were recalculated to slightly
# narrower widths (the workbook was re-saved by a newer Excel build with
# different font metrics). Reproduce the resulting widths as closely as
# this engine's width<->pixel rounding allows.
$ws.Columns("B").ColumnWidth = 16.0
$ws.Columns("C").ColumnWidth = 14.666667
$ws.Columns("D").ColumnWidth = 31.333333
$ws.Columns("E").ColumnWidth = 11.833333
$ws.Columns("F").ColumnWidth = 8.666667
$ws.Columns("G").ColumnWidth = 12.666667
$ws.Columns("H").ColumnWidth = 17.0
$ws.Columns("I").ColumnWidth = 16.0
$ws.Columns("J").ColumnWidth = 13.333333
